$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9 (G9=5487)
$ws.Range("H9").Value = 22
$ws.Range("I9").Value = 24.5
$ws.Range("J9").Value = 20.333334
$ws.Range("K9").Value = 24.5
$ws.Range("L9").Value = 20.333334
$ws.Range("M9").Value = 144.5
$ws.Range("N9").Value = -358.333334

# Row 18 (G18=5471)
$ws.Range("H18").Value = 552.25
$ws.Range("I18").Value = 552.25
$ws.Range("K18").Value = 552.25
$ws.Range("M18").Value = -268.25

# Row 19 (G19=7015)
$ws.Range("H19").Value = 1956.6666
$ws.Range("I19").Value = 3049.5
$ws.Range("J19").Value = 1410.25
$ws.Range("K19").Value = 3049.5
$ws.Range("L19").Value = 1410.25
$ws.Range("M19").Value = -2874.5
$ws.Range("N19").Value = -1760.25

# Row 43 (G43=5472)
$ws.Range("H43").Value = 2526.8572
$ws.Range("I43").Value = 1944
$ws.Range("J43").Value = 2760
$ws.Range("K43").Value = 1944
$ws.Range("L43").Value = 2760
$ws.Range("M43").Value = -1875
$ws.Range("N43").Value = -2898

# Row 80 (G80=12605)
$ws.Range("H80").Value = 29174.047
$ws.Range("I80").Value = 67001.11
$ws.Range("K80").Value = 201003.33
$ws.Range("M80").Value = -200005.33

# Row 83 (G83=12605)
$ws.Range("H83").Value = 29174.047
$ws.Range("I83").Value = 67001.11
$ws.Range("K83").Value = 603009.99
$ws.Range("M83").Value = -598017.99

# Row 129 (G129=36115)
$ws.Range("H129").Value = 715
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

# Row 132 (G132=44049)
$ws.Range("H132").Value = 3128.9722
$ws.Range("I132").Value = 2419.2058
$ws.Range("K132").Value = 7257.617400000001
$ws.Range("M132").Value = -4727.617400000001

# Row 137 (G137=44013)
$ws.Range("H137").Value = 2360.0715
$ws.Range("I137").Value = 1060
$ws.Range("J137").Value = 3335.125
$ws.Range("K137").Value = 3180
$ws.Range("L137").Value = 10005.375
$ws.Range("M137").Value = -630
$ws.Range("N137").Value = -15105.375

# Row 138 (G138=44169)
$ws.Range("H138").Value = 3500
$ws.Range("J138").Value = 3090.3704
$ws.Range("L138").Value = 9271.111199999999
$ws.Range("N138").Value = -19551.1112

$ws = $wb.Worksheets.Item("ARM")
# Row 122 (G122=36168)
$ws.Range("H122").Value = 2473.75
$ws.Range("I122").Value = 2022.5
$ws.Range("K122").Value = 6067.5
$ws.Range("M122").Value = -3617.5

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (G94=19939)
$ws.Range("H94").Value = 2365.4736
$ws.Range("I94").Value = 1917.1333
$ws.Range("K94").Value = 1917.1333
$ws.Range("M94").Value = -1466.1333

# Row 99 (G99=19943)
$ws.Range("H99").Value = 56449.21
$ws.Range("I99").Value = 86643.75
$ws.Range("K99").Value = 86643.75
$ws.Range("M99").Value = -85145.75

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (G22=5367)
$ws.Range("H22").Value = 1290.15
$ws.Range("I22").Value = 289.75
$ws.Range("J22").Value = 2790.75
$ws.Range("K22").Value = 289.75
$ws.Range("L22").Value = 2790.75
$ws.Range("M22").Value = 60.25
$ws.Range("N22").Value = -3490.75

# Row 31 (G31=44023)
$ws.Range("H31").Value = 2305.6296
$ws.Range("I31").Value = 1781.9474
$ws.Range("K31").Value = 1781.9474
$ws.Range("M31").Value = -1486.9474

# Row 34 (G34=44023)
$ws.Range("H34").Value = 2305.6296
$ws.Range("I34").Value = 1781.9474
$ws.Range("K34").Value = 1781.9474
$ws.Range("M34").Value = -1579.9474

# Row 86 (G86=12584)
$ws.Range("H86").Value = 47743.555
$ws.Range("I86").Value = 53117.75
$ws.Range("J86").Value = 4750
$ws.Range("K86").Value = 53117.75
$ws.Range("L86").Value = 4750
$ws.Range("M86").Value = -51994.75
$ws.Range("N86").Value = -6996

# Row 89 (G89=12584)
$ws.Range("H89").Value = 47743.555
$ws.Range("I89").Value = 53117.75
$ws.Range("J89").Value = 4750
$ws.Range("K89").Value = 265588.75
$ws.Range("L89").Value = 23750
$ws.Range("M89").Value = -259972.75
$ws.Range("N89").Value = -34982

# Row 94 (G94=32934)
$ws.Range("H94").Value = 3343
$ws.Range("I94").Value = 3236.5
$ws.Range("J94").Value = 3403.8572
$ws.Range("K94").Value = 3236.5
$ws.Range("L94").Value = 3403.8572
$ws.Range("M94").Value = -2785.5
$ws.Range("N94").Value = -4305.8572

# Row 99 (G99=36198)
$ws.Range("H99").Value = 2941.4443
$ws.Range("J99").Value = 3225.1428
$ws.Range("L99").Value = 3225.1428
$ws.Range("N99").Value = -6221.1428

# Row 107 (G107=27689)
$ws.Range("H107").Value = 389.56
$ws.Range("I107").Value = 377.5238
$ws.Range("K107").Value = 377.5238
$ws.Range("M107").Value = 1542.4762

# Row 126 (G126=36198)
$ws.Range("H126").Value = 2941.4443
$ws.Range("J126").Value = 3225.1428
$ws.Range("L126").Value = 9675.428400000001
$ws.Range("N126").Value = -14615.4284

$ws = $wb.Worksheets.Item("CUL")
# Row 80 (G80=12890)
$ws.Range("H80").Value = 5924.3335
$ws.Range("J80").Value = 5924.3335
$ws.Range("L80").Value = 17773.0005
$ws.Range("N80").Value = -19645.0005

# Row 83 (G83=12890)
$ws.Range("H83").Value = 5924.3335
$ws.Range("J83").Value = 5924.3335
$ws.Range("L83").Value = 53319.0015
$ws.Range("N83").Value = -62679.0015

# Row 106 (G106=19819)
$ws.Range("H106").Value = 13000
$ws.Range("J106").Value = 13000
$ws.Range("L106").Value = 39000
$ws.Range("N106").Value = -40892

# Row 109 (G109=27854)
$ws.Range("H109").Value = 2035.7142
$ws.Range("I109").Value = 1875
$ws.Range("K109").Value = 5625
$ws.Range("M109").Value = -4585

# Row 112 (G112=27855)
$ws.Range("H112").Value = 12466.667
$ws.Range("I112").Value = 2400
$ws.Range("J112").Value = 17500
$ws.Range("K112").Value = 7200
$ws.Range("L112").Value = 52500
$ws.Range("M112").Value = -6092
$ws.Range("N112").Value = -54716

# Row 113 (G113=27843)
$ws.Range("H113").Value = 926.75
$ws.Range("I113").Value = 920.5
$ws.Range("K113").Value = 2761.5
$ws.Range("M113").Value = -591.5

# Row 117 (G117=27870)
$ws.Range("H117").Value = 5332.6665
$ws.Range("I117").Value = 2499.5
$ws.Range("K117").Value = 7498.5
$ws.Range("M117").Value = -4056.5

# Row 129 (G129=36054)
$ws.Range("H129").Value = 783669.2
$ws.Range("I129").Value = 14307.75
$ws.Range("K129").Value = 42923.25
$ws.Range("M129").Value = -37923.25

# Row 131 (G131=36060)
$ws.Range("H131").Value = 2711444
$ws.Range("I131").Value = 18913.834
$ws.Range("J131").Value = 3232578.8
$ws.Range("K131").Value = 56741.50199999999
$ws.Range("L131").Value = 9697736.399999999
$ws.Range("M131").Value = -51701.50199999999
$ws.Range("N131").Value = -9707816.399999999

# Row 134 (G134=44074)
$ws.Range("H134").Value = 5491.2
$ws.Range("I134").Value = 5491.2
$ws.Range("K134").Value = 16473.6
$ws.Range("M134").Value = -11403.6

$ws = $wb.Worksheets.Item("GSM")
# Row 102 (G102=36169)
$ws.Range("H102").Value = 1150.5
$ws.Range("I102").Value = 853.38464
$ws.Range("K102").Value = 853.38464
$ws.Range("M102").Value = 768.61536

# Row 122 (G122=36182)
$ws.Range("H122").Value = 2717.182
$ws.Range("I122").Value = 2311.625
$ws.Range("J122").Value = 3798.6667
$ws.Range("K122").Value = 6934.875
$ws.Range("L122").Value = 11396.0001
$ws.Range("M122").Value = -4484.875
$ws.Range("N122").Value = -16296.0001

# Row 132 (G132=44008)
$ws.Range("H132").Value = 27580.46
$ws.Range("I132").Value = 32887.72
$ws.Range("K132").Value = 98663.16
$ws.Range("M132").Value = -96133.16

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (G7=36249)
$ws.Range("H7").Value = 4297.88
$ws.Range("I7").Value = 3399.3333
$ws.Range("K7").Value = 3399.3333
$ws.Range("M7").Value = -3287.3333

# Row 70 (G70=10811)
$ws.Range("H70").Value = 49999
$ws.Range("J70").Value = 49999
$ws.Range("L70").Value = 49999
$ws.Range("N70").Value = -50539

# Row 73 (G73=10811)
$ws.Range("H73").Value = 49999
$ws.Range("J73").Value = 49999
$ws.Range("L73").Value = 49999
$ws.Range("N73").Value = -51871

# Row 93 (G93=19993)
$ws.Range("H93").Value = 2042.8889
$ws.Range("I93").Value = 1848.5
$ws.Range("J93").Value = 2237.2778
$ws.Range("K93").Value = 1848.5
$ws.Range("L93").Value = 2237.2778
$ws.Range("M93").Value = -600.5
$ws.Range("N93").Value = -4733.2778

# Row 126 (G126=36249)
$ws.Range("H126").Value = 4297.88
$ws.Range("I126").Value = 3399.3333
$ws.Range("K126").Value = 10197.9999
$ws.Range("M126").Value = -7727.999899999999

$ws = $wb.Worksheets.Item("WVR")
# Row 62 (G62=12589)
$ws.Range("H62").Value = 110542.09
$ws.Range("J62").Value = 134107.11
$ws.Range("L62").Value = 134107.11
$ws.Range("N62").Value = -135355.11

# Row 65 (G65=12589)
$ws.Range("H65").Value = 110542.09
$ws.Range("J65").Value = 134107.11
$ws.Range("L65").Value = 670535.5499999999
$ws.Range("N65").Value = -676775.5499999999
